# Translation cleanup:
#  - Fix missing escape characters around the parentheses in the
#    "Global scale (marine or terrestrial)" / "Globale (marina o terrestre)" entry.
#  - Correct the spelling of "Palearctic" -> "Palaearctic" (English column only;
#    the Italian "Paleartica" is unchanged).
#  - Remove the stray trailing empty row left at the bottom of the table.
#  - Leave the active selection on B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order matters here only insofar as it keeps the shared-string table tidy;
# fix the spelling first, then the escaped-parentheses pair.
$ws.Range("A7").Value = "Palaearctic"

$ws.Range("A2").Value = "Global scale \(marine or terrestrial\)"
$ws.Range("B2").Value = "Globale \(marina o terrestre\)"

# Drop the empty trailing row (row 16) that had no data.
$ws.Rows("16").Delete()

# Match the saved selection/active cell.
$ws.Range("B2").Select()
